$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.901.42"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.501.83"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.64"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.24"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "2.942.57"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "58.826.31"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.75"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "2.498.77"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.04"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.82"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.51"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "0.0₃0761"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.10"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.37"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.34"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.04"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.57"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.798"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "281.52"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.602"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.90"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -5.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "128.59"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0925"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.22"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "1.749.24"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  -0.52%  "
